# Add data for 2022-07-29: refresh the "through" date from 2022-07-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab (workbook.xml <sheet name="...">)
$ws.Name = "Through 2022-07-21"

# Update the 2022-column header label (I1) to match the new "through" date
$ws.Range("I1").Value = "2022 (through 07-21)"

# Update the 2022 figures that changed: July (row 8) and the yearly Total (row 14)
$ws.Range("I8").Value = 118
$ws.Range("I14").Value = 924
